# Update the ROM ore cost block (rows 105-143) with refreshed model-run
# output values for columns A (open_pit1) through J (underground7).
# All 39 rows share the same flat value per column, matching the
# pre-edit data pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    16.3021472228724,
    50.09795341786098,
    91.96795783773388,
    133.1132357349561,
    142.5301089763919,
    151.9469822178277,
    161.3638554680526,
    170.7807287182775,
    180.1976019597133,
    187.6206087008836
)

$firstRow = 105
$lastRow = 143

for ($row = $firstRow; $row -le $lastRow; $row++) {
    for ($col = 1; $col -le $newValues.Length; $col++) {
        $ws.Cells.Item($row, $col).Value = $newValues[$col - 1]
    }
}
